$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Row 5 (the "1/14" row), Column 2 (time column): ---
# After the existing "9:05 - 9:50 AM" paragraph, add 14 empty paragraphs
# followed by a new paragraph containing "2:50 - 3:15 PM".
$cellTime = $t.Cell(5, 2)
$timePara = $cellTime.Range.Paragraphs.Item(1)
$timeRng = $timePara.Range
$timeRng.Collapse(0)

$breaks = ""
for ($i = 0; $i -lt 14; $i++) {
    $breaks = $breaks + "`r"
}
$timeRng.InsertAfter($breaks + "2:50 – 3:15 PM")

# --- Row 5 (the "1/14" row), Column 3 (accomplishments column): ---
# Append a new paragraph at the end of the cell.
$cellWork = $t.Cell(5, 3)
$workPara = $cellWork.Range.Paragraphs.Item(1)
$workRng = $workPara.Range
$workRng.Collapse(0)
$workRng.InsertAfter("`rTemporarily rewrote the KeyboardInput script of the Unite 2017 demo.")
